$wb = $excel.ActiveWorkbook

# --- Rename the last sheet ("Sheet1" -> "level constraint") and add the new data ---
$ws = $wb.Worksheets.Item(5)
$ws.Name = "level constraint"

# New "Time (s)" (G) and "Experimental" (I) values for each data row (2-16).
# J/K/L are existing formulas on the sheet; they recompute automatically once
# I (and H, already present) are populated.
$ws.Range("G2").Value = 0.01000952720642079956
$ws.Range("I2").Value = 3140

$ws.Range("G3").Value = 0.00400042533874511025
$ws.Range("I3").Value = 3166

$ws.Range("G4").Value = 0.00797820091247558073
$ws.Range("I4").Value = 3131

$ws.Range("G5").Value = 0.00997257232666014931
$ws.Range("I5").Value = 3149

$ws.Range("G6").Value = 0.00486922264099121007
$ws.Range("I6").Value = 3314

$ws.Range("G7").Value = 0.00399875640869140018
$ws.Range("I7").Value = 3282

$ws.Range("G8").Value = 0.00598406791687011025
$ws.Range("I8").Value = 3393

$ws.Range("G9").Value = 0.00797986984252928994
$ws.Range("I9").Value = 3467

$ws.Range("G10").Value = 0.01002717018127440018
$ws.Range("I10").Value = 3403

$ws.Range("G11").Value = 0.00799775123596191059
$ws.Range("I11").Value = 3692

$ws.Range("G12").Value = 0.00701737403869628039
$ws.Range("I12").Value = 4116

$ws.Range("G13").Value = 0.00599288940429686980
$ws.Range("I13").Value = 4287

$ws.Range("G14").Value = 0.01097321510314940018
$ws.Range("I14").Value = 4433

$ws.Range("G15").Value = 0.00602006912231444966
$ws.Range("I15").Value = 4341

$ws.Range("G16").Value = 0.00997495651245117014
$ws.Range("I16").Value = 4173

# Column widths for I and L (closest reachable approximation of the
# author's manual column resize: 13.85546875 and 10.5703125 respectively).
$ws.Columns.Item(9).ColumnWidth = 13
$ws.Columns.Item(12).ColumnWidth = 9.666666666666666

# Selection moves to M4 on this sheet.
$ws.Range("M4").Select() | Out-Null

# Keep the "level constraint" tab active/selected, matching the saved view.
$ws.Activate() | Out-Null
